$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.116.37'
$ws.Range("E2").Value = '  +1.48%  '
$ws.Range("D3").Value = '2.881.86'
$ws.Range("E3").Value = '  +3.88%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '352.69'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.79'
$ws.Range("E6").Value = '  +2.89%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.561'
$ws.Range("E7").Value = '  +2.38%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.17'
$ws.Range("E10").Value = '  +2.56%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0860'
$ws.Range("E11").Value = '  +3.38%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("E14").Value = '  +1.10%  '
$ws.Range("D15").Value = '3.338.17'
$ws.Range("E15").Value = '  +4.14%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '2.916.65'
$ws.Range("E16").Value = '  +4.98%  '
$ws.Range("B17").Value = 'Polygon'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.997'
$ws.Range("E17").Value = '  +8.14%  '
$ws.Range("D18").Value = '52.147.81'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.36'
$ws.Range("E19").Value = '  +8.22%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.71'
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.76'
$ws.Range("E21").Value = '  +2.93%  '
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.06'
$ws.Range("E23").Value = '  +1.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '270.20'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.79'
$ws.Range("E25").Value = '  +1.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.35'
$ws.Range("E26").Value = '  +2.22%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.56'
$ws.Range("E29").Value = '  +3.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.82'
$ws.Range("E30").Value = '  +5.11%  '
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '53.29'
$ws.Range("E33").Value = '  +2.84%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.92'
$ws.Range("E34").Value = '  +6.88%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0935'
$ws.Range("E35").Value = '  +11.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0458'
$ws.Range("E36").Value = '  +3.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  +7.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.67'
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("E40").Value = '  +3.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.66'
$ws.Range("E41").Value = '  +7.17%  '
$ws.Range("E42").Value = '  +3.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.75'
$ws.Range("E43").Value = '  +4.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '121.80'
$ws.Range("E44").Value = '  +1.53%  '
$ws.Range("E45").Value = '  +1.63%  '
$ws.Range("E46").Value = '  +6.41%  '
$ws.Range("D47").Value = '2.204.58'
$ws.Range("E47").Value = '  +3.60%  '
$ws.Range("E48").Value = '  +6.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.272'
$ws.Range("E49").Value = '  +20.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.954'
$ws.Range("E50").Value = '  +7.12%  '
$ws.Range("E51").Value = '  +1.62%  '
